$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (host, domain) pairs for rows 1-12 replacing the old IP/domain list.
# Column A is a formula ("sudo docker run ... -l "&B&"&& sleep 5;") that
# recalculates automatically once B is updated.
$data = @(
    @("195.161.52.80 ", "pfrf.ru"),
    @("193.148.44.189", "fss.ru"),
    @("193.148.44.186", "data.fss.ru "),
    @("193.148.44.187", "docs.fss.ru "),
    @("193.148.44.192", "sip.fss.ru"),
    @("193.148.44.195", "support.fss.ru"),
    @("193.148.44.196", "map.fss.ru"),
    @("193.148.44.202", "ftp.fss.ru"),
    @("193.148.44.208", "av.fss.ru "),
    @("193.148.44.217", "fw.fss.ru "),
    @("193.148.44.218", "proxy.fss.ru"),
    @("193.148.44.29 ", "beta.fss.ru")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Range("B$row").Value = $data[$i][0]
    $ws.Range("C$row").Value = $data[$i][1]
}

# Rows 13-69 previously held additional (formula/B/C) entries; they are now
# emptied out entirely (cells keep their existing style, just no content).
for ($row = 13; $row -le 69; $row++) {
    $ws.Range("A$row`:C$row").ClearContents()
}
